# Clean up the counters summary data: zero out the COMPLETENESSMANDATORY
# counts/scores (columns B and C) for the attributes that previously had
# non-zero values. This reflects the reworked data-processing logic where
# the operation driving these counters is now fully controllable and the
# metadata is generated fresh per dataset (starting at zero).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where column B (COMPLETENESSMANDATORY) must become 0
$rowsB = @(6, 7, 8, 9, 20, 22, 23, 24, 25, 28, 32, 33, 34, 36)
foreach ($r in $rowsB) {
    $ws.Cells.Item($r, 2).Value = 0
}

# Rows where column C (COMPLETENESSMANDATORY SCORE) must also become 0
$rowsC = @(8, 9, 22, 25, 28, 34, 36)
foreach ($r in $rowsC) {
    $ws.Cells.Item($r, 3).Value = 0
}
